# Meerbeek data request - annual reports, harvest spreadsheet.
# Fills in the 2017-2019 lake survey data (catch/hr, mean weight, population
# estimate, biomass, acres) and the related footnotes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A to fit the longer lake names (e.g. "North Twin")
$ws.Columns.Item(1).ColumnWidth = 13.85546875

# ---- 2017 ----
$ws.Range("A2").Value = 'Blue'
$ws.Range("B2").Value = 2017
$ws.Range("E2").Value = 60835
$ws.Range("F2").Value = 448.3903345724907
$ws.Range("F2").NumberFormat = "0.0"
$ws.Range("G2").Value = 269
$ws.Range("A3").Value = 'Center'
$ws.Range("B3").Value = 2017
$ws.Range("E3").Value = 8675
$ws.Range("F3").Value = 263.32727272727271
$ws.Range("F3").NumberFormat = "0.0"
$ws.Range("G3").Value = 220
$ws.Range("A4").Value = 'North Twin'
$ws.Range("B4").Value = 2017
$ws.Range("E4").Value = 2057
$ws.Range("F4").Value = 19.869757174392937
$ws.Range("F4").NumberFormat = "0.0"
$ws.Range("G4").Value = 453
$ws.Range("A5").Value = 'South Twin'
$ws.Range("B5").Value = 2017
$ws.Range("E5").Value = 23204
$ws.Range("F5").Value = 68.644999999999996
$ws.Range("F5").NumberFormat = "0.0"
$ws.Range("G5").Value = 600

# ---- 2018 ----
$ws.Range("A6").Value = 'Blue'
$ws.Range("B6").Value = 2018
$ws.Range("C6").Value = 75.900000000000006
$ws.Range("E6").Value = 72140
$ws.Range("F6").Value = 543
$ws.Range("F6").NumberFormat = "0.0"
$ws.Range("G6").Value = 269
$ws.Range("A7").Value = 'Center'
$ws.Range("B7").Value = 2018
$ws.Range("C7").Value = 5.47
$ws.Range("E7").Value = 6466
$ws.Range("F7").Value = 194
$ws.Range("F7").NumberFormat = "0.0"
$ws.Range("G7").Value = 220
$ws.Range("A8").Value = 'Five Island'
$ws.Range("B8").Value = 2018
$ws.Range("C8").Value = 10.8
$ws.Range("E8").Value = 25798
$ws.Range("F8").Value = 187
$ws.Range("F8").NumberFormat = "0.0"
$ws.Range("G8").Value = 973
$ws.Range("A9").Value = 'North Twin'
$ws.Range("B9").Value = 2018
$ws.Range("C9").Value = 41.3
$ws.Range("E9").Value = 3816
$ws.Range("F9").Value = 34
$ws.Range("F9").NumberFormat = "0.0"
$ws.Range("G9").Value = 453
$ws.Range("A10").Value = 'South Twin'
$ws.Range("B10").Value = 2018
$ws.Range("C10").Value = 105
$ws.Range("E10").Value = 20661
$ws.Range("F10").Value = 62
$ws.Range("F10").NumberFormat = "0.0"
$ws.Range("G10").Value = 600
$ws.Range("A11").Value = 'Silver'
$ws.Range("B11").Value = 2018
$ws.Range("C11").Value = 14.1
$ws.Range("E11").Value = 9755
$ws.Range("F11").Value = 96
$ws.Range("F11").NumberFormat = "0.0"
$ws.Range("G11").Value = 1041
$ws.Range("A12").Value = 'Storm'
$ws.Range("B12").Value = 2018
$ws.Range("C12").Value = 12.9
$ws.Range("E12").Value = 9251
$ws.Range("F12").Value = 11
$ws.Range("F12").NumberFormat = "0.0"
$ws.Range("G12").Value = 3097

# ---- 2019 ----
$ws.Range("A13").Value = 'Blue'
$ws.Range("B13").Value = 2019
$ws.Range("C13").Value = 191
$ws.Range("E13").Value = 25661
$ws.Range("F13").Value = 192
$ws.Range("F13").NumberFormat = "0.0"
$ws.Range("G13").Value = 269
$ws.Range("A14").Value = 'Center'
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 6.1
$ws.Range("E14").Value = 1451
$ws.Range("F14").Value = 48
$ws.Range("F14").NumberFormat = "0.0"
$ws.Range("G14").Value = 220
$ws.Range("A15").Value = 'Five Island'
$ws.Range("B15").Value = 2019
$ws.Range("C15").Value = 9.81
$ws.Range("E15").Value = 19738
$ws.Range("F15").Value = 160.19999999999999
$ws.Range("F15").NumberFormat = "0.0"
$ws.Range("G15").Value = 973
$ws.Range("A16").Value = 'North Twin'
$ws.Range("B16").Value = 2019
$ws.Range("C16").Value = 30.4
$ws.Range("E16").Value = 2487
$ws.Range("F16").Value = 34.5
$ws.Range("F16").NumberFormat = "0.0"
$ws.Range("G16").Value = 453
$ws.Range("A17").Value = 'South Twin'
$ws.Range("B17").Value = 2019
$ws.Range("C17").Value = 195
$ws.Range("E17").Value = 14896
$ws.Range("F17").Value = 45.3
$ws.Range("F17").NumberFormat = "0.0"
$ws.Range("G17").Value = 600
$ws.Range("A18").Value = 'Silver'
$ws.Range("B18").Value = 2019
$ws.Range("C18").Value = 23.1
$ws.Range("E18").Value = 9174
$ws.Range("F18").Value = 93.8
$ws.Range("F18").NumberFormat = "0.0"
$ws.Range("G18").Value = 1041
$ws.Range("A19").Value = 'Storm'
$ws.Range("B19").Value = 2019
$ws.Range("C19").Value = 20.3
$ws.Range("E19").Value = 15467
$ws.Range("F19").Value = 20.6
$ws.Range("F19").NumberFormat = "0.0"
$ws.Range("G19").Value = 3097

# ---- Footnotes (rows 22-25) ----
# Rows 22/24 reuse the existing footnote text (previously at A10/A11).
# Rows 23/25 are new highlighted notes below each footnote.
$ws.Range("A22").Value = '*CAP catch per hour should only be from bimonthly EF samples collected from August to October (combined)'
$ws.Range("B23").Value = 'No fall standard runs in 2017'
$ws.Range("B23").Interior.Color = 65535
$ws.Range("A24").Value = '** CAP mean wt should only be from those fish captured during bimonthly EF samples from August to October (combined)'
$ws.Range("B25").Value = "Calculated from each lake's spring l-w regression, applied to observed fall catch lengths (std runs), then averaged"
$ws.Range("B25").Interior.Color = 65535

# Update selection to match final cursor position
$ws.Range("B26").Select() | Out-Null
